$wb = $excel.ActiveWorkbook

# Rename sheets (tab names) - order matches workbook.xml sheet order
$wb.Worksheets.Item(1).Name = "GNG_TO-16509961345051398"
$wb.Worksheets.Item(2).Name = "NB_TO-16509961368097637"
$wb.Worksheets.Item(3).Name = "RS_TO-16509961368097637"
$wb.Worksheets.Item(4).Name = "TOL_TO-16509961368657374"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16509961369457672"

# Sheet 1: GNG_TO
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16509961344731069.csv"
$ws1.Range("B3").Value = "GNG_stims-16509961344891498.csv"
$ws1.Range("B4").Value = "go_stims-16509961344891498.csv"
$ws1.Range("B5").Value = "GNG_stims-16509961345051398.csv"

# Sheet 2: NB_TO
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "OB-16509961353214865.csv"
$ws2.Range("B3").Value = "ZB-match_3-16509961349934807.csv"
$ws2.Range("B4").Value = "ZB-match_8-16509961346411076.csv"
$ws2.Range("B5").Value = "TB-16509961365457678.csv"
$ws2.Range("B6").Value = "ZB-match_3-16509961348414824.csv"
$ws2.Range("B7").Value = "TB-16509961358817327.csv"
$ws2.Range("B8").Value = "OB-1650996135737484.csv"
$ws2.Range("B9").Value = "OB-16509961357055168.csv"
$ws2.Range("B10").Value = "TB-16509961367777681.csv"

# Sheet 3: RS_TO
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "eyes open"
$ws3.Range("B3").Value = "eyes closed"

# Sheet 4: TOL_TO
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16509961368337662.csv"
$ws4.Range("B3").Value = "ZM_stims-16509961368097637.csv"
$ws4.Range("B4").Value = "MM_stims-1650996136849732.csv"
$ws4.Range("B5").Value = "ZM_stims-16509961368337662.csv"
$ws4.Range("B6").Value = "MM_stims-16509961368657374.csv"
$ws4.Range("B7").Value = "ZM_stims-1650996136849732.csv"

# Sheet 5: vSAT_TO
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-16509961368977678.csv"
$ws5.Range("B3").Value = "vSAT_stims-16509961369297676.csv"
$ws5.Range("B4").Value = "SAT_stims-16509961368657374.csv"
$ws5.Range("B5").Value = "vSAT_stims-16509961369137378.csv"
